$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# --- New "category" column, inserted right after "property_category" (H) ---
# This shifts the former I/J/K (date/legislator_name/legislator_id) one
# column to the right, becoming J/K/L, and carries their formatting along.
$ws.Columns.Item(9).Insert()

$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(2, 9).Value = "normal"
$ws.Cells.Item(3, 9).Value = "normal"
$ws.Cells.Item(4, 9).Value = "normal"
$ws.Cells.Item(5, 9).Value = "normal"

# --- Two new trailing columns: "source_file" (M) and "index" (N) ---
# Clone the formatting of the adjacent "legislator_id" column (L) onto the
# new columns before filling in their values, so the header row keeps the
# same bold/bordered look as the rest of row 1.
$ws.Range("L1:L5").Copy()
$ws.Range("M1:M5").PasteSpecial(-4122)
$ws.Range("L1:L5").Copy()
$ws.Range("N1:N5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(2, 13).Value = "tmp7091"
$ws.Cells.Item(3, 13).Value = "tmp7091"
$ws.Cells.Item(4, 13).Value = "tmp7091"
$ws.Cells.Item(5, 13).Value = "tmp7091"

$ws.Cells.Item(1, 14).Value = "index"
$ws.Cells.Item(2, 14).Value = 58
$ws.Cells.Item(3, 14).Value = 59
$ws.Cells.Item(4, 14).Value = 60
$ws.Cells.Item(5, 14).Value = 61
